$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '22.409.15'
$ws.Range("E2").Value = '  +0.21%  '
$ws.Range("D3").Value = '1.571.05'
$ws.Range("E3").Value = '  +0.06%  '
$ws.Range("E4").Value = '  +0.22%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '1.002'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +0.19%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '290.16'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +0.00%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3744'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -1.47%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '49.58'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +0.11%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3369'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -1.79%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07471'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -2.84%  '
$ws.Range("E11").Value = '  -3.24%  '
$ws.Range("E12").Value = '  +0.27%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '21.02'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -1.85%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.918'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -1.99%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.869'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -1.22%  '
$ws.Range("D16").Value = '1.570.67'
$ws.Range("E16").Value = '  -0.52%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001116'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -2.22%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '89.35'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -1.34%  '
$ws.Range("E19").Value = '  -0.48%  '
$ws.Range("E20").Value = '  +0.24%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.166'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -1.57%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '16.16'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -3.13%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.85'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -1.27%  '
$ws.Range("D24").Value = '22.408.80'
$ws.Range("E24").Value = '  +0.11%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.367'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -0.92%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.552'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -8.67%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.03'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -1.21%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '147.01'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +1.54%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '4.992'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -0.62%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '124.69'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -1.03%  '
$ws.Range("D31").Value = '1.745.50'
$ws.Range("E31").Value = '  -0.26%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.9931'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -2.72%  '
$ws.Range("E33").Value = '  -3.28%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.916'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -5.04%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '9.721'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -4.25%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.08423'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -1.67%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.379'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +4.29%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02449'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -3.72%  '
$ws.Range("B39").Value = 'Hedera'
$ws.Range("C39").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.06469'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +1.18%  '
$ws.Range("B40").Value = 'Algorand'
$ws.Range("C40").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.2248'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -3.50%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.376'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -3.59%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '11.30'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -3.60%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.6206'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -3.33%  '
$ws.Range("B44").Value = 'EnergySwap'
$ws.Range("C44").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '14.09'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -1.03%  '
$ws.Range("B45").Value = 'Frax'
$ws.Range("C45").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.002'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +0.25%  '
$ws.Range("E46").Value = '  +1.34%  '
$ws.Range("E47").Value = '  -3.73%  '
$ws.Range("E48").Value = '  -2.22%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '125.55'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +0.87%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.226'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -8.02%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.07295'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +0.08%  '
